$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quantities (column B) to reflect new SAP-driven values
$ws.Range("B2").Value = 4500
$ws.Range("B3").Value = 1000
$ws.Range("B4").Value = 1000

# Remove the now-duplicate last row (row 5), shrinking the used range to A1:C4
$ws.Rows.Item(5).Delete()

# Restore the selection to where the user last left it
$ws.Range("C6").Select()
